$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell A3 to the new parameter name (creates a new shared string)
$ws.Range("A3").Value = "cb_distances_calc"

# Move the active selection to A3 (matches the saved selection state in the diff)
$ws.Range("A3").Select()
